$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update matricule (employee ID) values in column B for the leave rows.
$ws.Range("B2").Value = 30046391
$ws.Range("B3").Value = 30046179
$ws.Range("B4").Value = 30046112
$ws.Range("B5").Value = 30045885

# Row height tweak that came along with the re-generated rows.
$ws.Rows.Item(1).RowHeight = 19.5
$ws.Rows.Item(5).RowHeight = 19.5
